$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values: A1 becomes the email, B1 becomes "xyz",
# C1/D1/E1 keep their original text (resume path / status / timestamp).
$ws.Range("B1").Value = "xyz"
$ws.Range("A1").Value = "xyz.qa@gmail.com"
$ws.Range("C1").Value = "D:\\hari vardhan\\Resume.doc"
$ws.Range("D1").Value = "Resume Uploaded Successfully"
$ws.Range("E1").Value = "Tue May 31 09:18:27 IST 2016"

# Move the active selection to A8 (reflects where the cursor was left
# when the workbook was last saved).
$ws.Range("A8").Select()
